$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.765.94"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "1.756.86"
$ws.Range("E3").Value = "  -3.14%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'320.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'0.4226"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("D8").Value = "'0.3618"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("D9").Value = "'42.37"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.60%  "
$ws.Range("D10").Value = "'0.07440"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("D11").Value = "'1.085"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.87%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "'20.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.45%  "
$ws.Range("D14").Value = "'6.045"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("D15").Value = "'7.273"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("D16").Value = "1.758.83"
$ws.Range("E16").Value = "  -3.62%  "
$ws.Range("D17").Value = "'90.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.64%  "
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").Value = "'0.06346"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'17.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "'5.917"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.73%  "
$ws.Range("D23").Value = "27.787.84"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").Value = "'11.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.50%  "
$ws.Range("D25").Value = "'2.104"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'157.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").Value = "'20.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.32%  "
$ws.Range("D28").Value = "1.961.73"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("D29").Value = "'2.128"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.89%  "
$ws.Range("D30").Value = "'123.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("D31").Value = "'1.115"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.47%  "
$ws.Range("D32").Value = "'3.685"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'5.540"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.23%  "
$ws.Range("D34").Value = "'0.08815"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.74%  "
$ws.Range("D35").Value = "'12.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.00%  "
$ws.Range("D36").Value = "'0.02286"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.86%  "
$ws.Range("D37").Value = "'0.06022"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.42%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("D39").Value = "'0.6292"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").Value = "'4.931"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.28%  "
$ws.Range("D41").Value = "'1.178"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").Value = "'0.9989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'1.396"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("D44").Value = "'7.831"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'13.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("D46").Value = "'0.5856"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").Value = "'122.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.92%  "
$ws.Range("D49").Value = "'1.972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.27%  "
$ws.Range("D50").Value = "'1.175"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "'0.06807"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.58%  "
